$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: leave blank but apply yellow fill across A:F (matches new style s="4": fillId=2, no border)
$row22 = $ws.Range("A22:F22")
$row22.Interior.Color = 65535

# Row 23: new experiment data
$ws.Range("A23").Value = "Exp 26"
$ws.Range("B23").Value = 0.1
$ws.Range("C23").Value = 10
$ws.Range("D23").Value = "Local"
$ws.Range("E23").Value = -1

# Apply same style as rest of data rows (s="3") to the new row's populated cells
$ws.Range("A23:E23").HorizontalAlignment = -4108

# Update view: scroll so topLeftCell = A5 and selection = G23
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("G23").Select()
